# pypsa2smspp/data/smspp_parameters.xlsx -- "First version capacity expansion uc"
#
# Adds InvestmentCost / MaxCapacityDesign parameter rows to the
# IntermittentUnitBlock sheet and Battery/Converter investment-cost +
# max-capacity-design parameter rows to the BatteryUnitBlock sheet, then
# moves the active tab/selection from BatteryUnitBlock back to
# IntermittentUnitBlock.

$wb = $excel.ActiveWorkbook

$xlCenter = -4108

# ---------------------------------------------------------------------
# BatteryUnitBlock: 4 new rows (Battery/Converter InvestmentCost +
# Battery/Converter MaxCapacityDesign) appended after the existing data.
# ---------------------------------------------------------------------
$wsBattery = $wb.Worksheets.Item("BatteryUnitBlock")

$wsBattery.Range("A17").Value2 = "BatteryInvestmentCost"
$wsBattery.Range("A18").Value2 = "ConverterInvestmentCost"
$wsBattery.Range("A19").Value2 = "BatteryMaxCapacityDesign"
$wsBattery.Range("A20").Value2 = "ConverterMaxCapacityDesign"

$wsBattery.Range("B17").Value2 = "Battery investment cost"
$wsBattery.Range("B18").Value2 = "Converter investment cost"
$wsBattery.Range("B19").Value2 = "Battery max capacity design"
$wsBattery.Range("B20").Value2 = "Converter max capacity design"

for ($r = 17; $r -le 20; $r++) {
    $wsBattery.Range("C$r").Value2 = "Parameter"
    $wsBattery.Range("D$r").Value2 = "float"
    $wsBattery.Range("E$r").Value2 = 1
    # Copy an existing "True" text cell so the new cell keeps the shared
    # string ("True") rather than becoming a native boolean.
    [void]$wsBattery.Range("F2").Copy()
    [void]$wsBattery.Range("F$r").PasteSpecial()
}

$wsBattery.Range("A17:F20").HorizontalAlignment = $xlCenter

# ---------------------------------------------------------------------
# IntermittentUnitBlock: 2 new rows (InvestmentCost + MaxCapacityDesign)
# appended after the existing data.
# ---------------------------------------------------------------------
$wsIntermittent = $wb.Worksheets.Item("IntermittentUnitBlock")

$wsIntermittent.Range("A8").Value2 = "InvestmentCost"
$wsIntermittent.Range("A9").Value2 = "MaxCapacityDesign"

$wsIntermittent.Range("B8").Value2 = "Investment cost"
$wsIntermittent.Range("B9").Value2 = "Max size"

for ($r = 8; $r -le 9; $r++) {
    $wsIntermittent.Range("C$r").Value2 = "Parameter"
    $wsIntermittent.Range("D$r").Value2 = "float"
    $wsIntermittent.Range("E$r").Value2 = 1
    [void]$wsIntermittent.Range("F2").Copy()
    [void]$wsIntermittent.Range("F$r").PasteSpecial()
}

$wsIntermittent.Range("A8:F9").HorizontalAlignment = $xlCenter

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Move the active tab / selection from BatteryUnitBlock to
# IntermittentUnitBlock (selecting the full A:F columns on each).
# ---------------------------------------------------------------------
[void]$wsBattery.Range("A1:F1048576").Select()
[void]$wsIntermittent.Activate()
[void]$wsIntermittent.Range("A1:F1048576").Select()
